$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'247.65"
$ws.Range("D3").Value2 = "'21.90"
$ws.Range("D4").Value2 = "'5.379"
$ws.Range("D5").Value2 = "'0.05637"
$ws.Range("D6").Value2 = "'3.428"
$ws.Range("D8").Value2 = "'0.8192"
$ws.Range("D9").Value2 = "'0.9341"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value2 = "'0.1446"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value2 = "'0.07481"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value2 = "'0.03243"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "ProBitToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D13").Value2 = "'0.1319"
$ws.Range("E13").Value = "12ProBitTokenPROB"
$ws.Range("D14").Value2 = "'0.03106"
$ws.Range("D15").Value2 = "'0.09318"
$ws.Range("D16").Value2 = "'3.570"
$ws.Range("D17").Value2 = "'0.001597"
$ws.Range("D18").Value2 = "'0.04726"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value2 = "'0.0005781"
$ws.Range("E19").Value = "18OneONE"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value2 = "'0.006368"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value2 = "'0.005057"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value2 = "'0.001032"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value2 = "'0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value2 = "'3.754"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value2 = "'2.188"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value2 = "'0.3309"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("D28").Value2 = "'0.0003000"
$ws.Range("D40").Value2 = "'0.03956"
$ws.Range("D41").Value2 = "'0.007001"
$ws.Range("D42").Value2 = "'0.1066"
$ws.Range("D43").Value2 = "'0.003021"
$ws.Range("D44").Value2 = "'0.008527"
$ws.Range("D45").Value2 = "'0.00005581"
$ws.Range("D47").Value2 = "'0.0005501"
$ws.Range("D48").Value2 = "'0.7801"
$ws.Range("D49").Value2 = "'0.1777"
$ws.Range("D50").Value2 = "'0.00002100"
